# The workbook previously used a volatile RANDBETWEEN-driven CHOOSE() formula
# in column F (F2 standalone, F3:F11 as a shared formula) to randomly assign
# one of {1,2,4}. This upload freezes that column to its last-calculated
# values, replacing the formulas with plain literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New frozen values for F2:F11 (row -> value)
$ws.Range("F2").Value  = 1
$ws.Range("F3").Value  = 1
$ws.Range("F4").Value  = 1
$ws.Range("F5").Value  = 4
$ws.Range("F6").Value  = 4
$ws.Range("F7").Value  = 4
$ws.Range("F8").Value  = 2
$ws.Range("F9").Value  = 2
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 2

# Leave the cursor on F10, matching the saved selection in the file.
$ws.Range("F10").Select()
